# AutoDS Proto enchaint auto optim + analyses => futur MCDSOptanalyser
#
# Adds a new "Params3_expl" sheet (MultiOpt-capable params, used to prototype
# the future MCDSOptanalyser), and tags the two existing Params sheets with
# the new "MultiOpt" column header.

$wb = $excel.ActiveWorkbook

$wsEchant1 = $wb.Worksheets.Item("Echant1_impl")
$wsEchant2 = $wb.Worksheets.Item("Echant2_impl")
$wsModl    = $wb.Worksheets.Item("Modl_impl")
$wsParams1 = $wb.Worksheets.Item("Params1_expl")
$wsParams2 = $wb.Worksheets.Item("Params2_expl")

# ---------------------------------------------------------------------------
# 1) Insert the new "Params3_expl" sheet right after "Params2_expl"
# ---------------------------------------------------------------------------
$wsParams3 = $wb.Worksheets.Add($null, $wsParams2)
$wsParams3.Name = "Params3_expl"

# ---------------------------------------------------------------------------
# 2) Populate Params3_expl (new MultiOpt-aware params table).
#    Rows are filled from the bottom up (row 7 -> row 1) so that the new
#    shared-string entries land in the same order as the target workbook.
# ---------------------------------------------------------------------------

# Row 7 : Luscinia megarhynchos / b / m / 5mn / 50 / auto / auto / times(2)
$wsParams3.Range("A7").Value = "Luscinia megarhynchos"
$wsParams3.Range("B7").Value = "b"
$wsParams3.Range("C7").Value = "m"
$wsParams3.Range("D7").Value = "5mn"
$wsParams3.Range("E7").Value = 50
$wsParams3.Range("F7").Value = "auto"
$wsParams3.Range("G7").Value = "auto"
$wsParams3.Range("H7").Value = "times(2)"

# Row 6 : Turdus merula / a+b / m / 10mn / auto / 200
$wsParams3.Range("A6").Value = "Turdus merula"
$wsParams3.Range("B6").Value = "a+b"
$wsParams3.Range("C6").Value = "m"
$wsParams3.Range("D6").Value = "10mn"
$wsParams3.Range("E6").Value = "auto"
$wsParams3.Range("F6").Value = 200

# Row 5 : Sylvia atricapilla / a+b / m / 5mn / auto / auto / 10
$wsParams3.Range("A5").Value = "Sylvia atricapilla"
$wsParams3.Range("B5").Value = "a+b"
$wsParams3.Range("C5").Value = "m"
$wsParams3.Range("D5").Value = "5mn"
$wsParams3.Range("E5").Value = "auto"
$wsParams3.Range("F5").Value = "auto"
$wsParams3.Range("G5").Value = 10

# Row 4 : Luscinia megarhynchos / b / m / 5mn / auto / auto
$wsParams3.Range("A4").Value = "Luscinia megarhynchos"
$wsParams3.Range("B4").Value = "b"
$wsParams3.Range("C4").Value = "m"
$wsParams3.Range("D4").Value = "5mn"
$wsParams3.Range("F4").Value = "auto"
$wsParams3.Range("G4").Value = "auto"

# Row 3 : Turdus merula / a+b / m / 10mn / auto / times(3, b=2)
$wsParams3.Range("A3").Value = "Turdus merula"
$wsParams3.Range("B3").Value = "a+b"
$wsParams3.Range("C3").Value = "m"
$wsParams3.Range("D3").Value = "10mn"
$wsParams3.Range("E3").Value = "auto"
$wsParams3.Range("H3").Value = "times(3, b=2)"

# Row 2 : Sylvia atricapilla / a+b / m / 5mn / auto / auto / times(4)
$wsParams3.Range("A2").Value = "Sylvia atricapilla"
$wsParams3.Range("B2").Value = "a+b"
$wsParams3.Range("C2").Value = "m"
$wsParams3.Range("D2").Value = "5mn"
$wsParams3.Range("E2").Value = "auto"
$wsParams3.Range("F2").Value = "auto"
$wsParams3.Range("H2").Value = "times(4)"

# Row 1 : header, added last
$wsParams3.Range("A1").Value = "Espèce"
$wsParams3.Range("B1").Value = "Passage"
$wsParams3.Range("C1").Value = "Adulte"
$wsParams3.Range("D1").Value = "Durée"
$wsParams3.Range("E1").Value = "TrGche"
$wsParams3.Range("F1").Value = "TrDrte"
$wsParams3.Range("G1").Value = "NbTrchMod"
$wsParams3.Range("H1").Value = "MultiOpt"

# Give the new header row the same bold/centered/bordered look as the other
# sheets' header rows (style index already used by A1:G1 on Params1/2_expl).
$wsParams2.Range("A1:G1").Copy() | Out-Null
$wsParams3.Range("A1:G1").PasteSpecial(-4122) | Out-Null

$wsParams2.Range("G1").Copy() | Out-Null
$wsParams3.Range("H1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3) Tag Params1_expl / Params2_expl with the new "MultiOpt" column header
# ---------------------------------------------------------------------------
$wsParams1.Range("G1").Copy() | Out-Null
$wsParams1.Range("H1").PasteSpecial(-4122) | Out-Null
$wsParams1.Range("H1").Value = "MultiOpt"

$wsParams2.Range("G1").Copy() | Out-Null
$wsParams2.Range("H1").PasteSpecial(-4122) | Out-Null
$wsParams2.Range("H1").Value = "MultiOpt"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Restore / update the per-sheet selections, then land on Params1_expl as
#    the active tab (matching the saved workbook view).
# ---------------------------------------------------------------------------
$wsEchant1.Activate()
$wsEchant1.Range("D27").Select() | Out-Null

$wsModl.Activate()
$wsModl.Range("I21").Select() | Out-Null

$wsParams2.Activate()
$wsParams2.Range("H1").Select() | Out-Null

$wsParams3.Activate()
$wsParams3.Range("H1").Select() | Out-Null

$wsParams1.Activate()
$wsParams1.Range("I9").Select() | Out-Null
